# Updates cryptocurrency price (column D) and 1-hour volume change (column E)
# values on the active worksheet to reflect the latest GitHub Actions-scraped
# snapshot, per the commit "Updated cryptos list ... with GitHub Actions".
#
# Price cells (column D) are plain text in the source data (e.g. "27.282.05"
# uses dot-grouping, not a valid number) so we force the Text number format
# before writing the value -- otherwise Excel's COM layer will silently
# coerce numeric-looking strings (e.g. "324.41", "0.00001026") into real
# numbers. The format is reset back to the default "Normal" style right
# after the write so the cell's style stays untouched (these cells carry no
# explicit style in the source file).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("D2").NumberFormat = "@"
    $ws.Range("D2").Value = '27.282.05'
    $ws.Range("D2").Style = "Normal"
    $ws.Range("E2").Value = '  -3.07%  '
    $ws.Range("D3").NumberFormat = "@"
    $ws.Range("D3").Value = '1.853.15'
    $ws.Range("D3").Style = "Normal"
    $ws.Range("E3").Value = '  -3.79%  '
    $ws.Range("E4").Value = '  +0.08%  '
    $ws.Range("D5").NumberFormat = "@"
    $ws.Range("D5").Value = '324.41'
    $ws.Range("D5").Style = "Normal"
    $ws.Range("E5").Value = '  -1.82%  '
    $ws.Range("E6").Value = '  +0.09%  '
    $ws.Range("D7").NumberFormat = "@"
    $ws.Range("D7").Value = '0.4542'
    $ws.Range("D7").Style = "Normal"
    $ws.Range("E7").Value = '  -3.85%  '
    $ws.Range("D8").NumberFormat = "@"
    $ws.Range("D8").Value = '0.3878'
    $ws.Range("D8").Style = "Normal"
    $ws.Range("E8").Value = '  -4.39%  '
    $ws.Range("D9").NumberFormat = "@"
    $ws.Range("D9").Value = '48.47'
    $ws.Range("D9").Style = "Normal"
    $ws.Range("E9").Value = '  -8.76%  '
    $ws.Range("D10").NumberFormat = "@"
    $ws.Range("D10").Value = '0.07905'
    $ws.Range("D10").Style = "Normal"
    $ws.Range("E10").Value = '  -6.17%  '
    $ws.Range("D11").NumberFormat = "@"
    $ws.Range("D11").Value = '1.013'
    $ws.Range("D11").Style = "Normal"
    $ws.Range("E11").Value = '  -3.43%  '
    $ws.Range("D12").NumberFormat = "@"
    $ws.Range("D12").Value = '21.36'
    $ws.Range("D12").Style = "Normal"
    $ws.Range("E12").Value = '  -4.21%  '
    $ws.Range("D13").NumberFormat = "@"
    $ws.Range("D13").Value = '1.847.80'
    $ws.Range("D13").Style = "Normal"
    $ws.Range("E13").Value = '  -3.17%  '
    $ws.Range("D14").NumberFormat = "@"
    $ws.Range("D14").Value = '5.899'
    $ws.Range("D14").Style = "Normal"
    $ws.Range("E14").Value = '  -3.30%  '
    $ws.Range("D15").NumberFormat = "@"
    $ws.Range("D15").Value = '7.138'
    $ws.Range("D15").Style = "Normal"
    $ws.Range("E15").Value = '  -4.99%  '
    $ws.Range("D16").NumberFormat = "@"
    $ws.Range("D16").Value = '1.003'
    $ws.Range("D16").Style = "Normal"
    $ws.Range("E16").Value = '  +0.24%  '
    $ws.Range("D17").NumberFormat = "@"
    $ws.Range("D17").Value = '0.06609'
    $ws.Range("D17").Style = "Normal"
    $ws.Range("E17").Value = '  +0.42%  '
    $ws.Range("D18").NumberFormat = "@"
    $ws.Range("D18").Value = '85.74'
    $ws.Range("D18").Style = "Normal"
    $ws.Range("E18").Value = '  -5.23%  '
    $ws.Range("D19").NumberFormat = "@"
    $ws.Range("D19").Value = '0.00001026'
    $ws.Range("D19").Style = "Normal"
    $ws.Range("E19").Value = '  -3.81%  '
    $ws.Range("D20").NumberFormat = "@"
    $ws.Range("D20").Value = '17.15'
    $ws.Range("D20").Style = "Normal"
    $ws.Range("E20").Value = '  -5.23%  '
    $ws.Range("E21").Value = '  +0.14%  '
    $ws.Range("D22").NumberFormat = "@"
    $ws.Range("D22").Value = '5.482'
    $ws.Range("D22").Style = "Normal"
    $ws.Range("E22").Value = '  -4.65%  '
    $ws.Range("D23").NumberFormat = "@"
    $ws.Range("D23").Value = '27.286.99'
    $ws.Range("D23").Style = "Normal"
    $ws.Range("E23").Value = '  -3.03%  '
    $ws.Range("E24").Value = '  -4.62%  '
    $ws.Range("E25").Value = '  +0.67%  '
    $ws.Range("D26").NumberFormat = "@"
    $ws.Range("D26").Value = '2.060.47'
    $ws.Range("D26").Style = "Normal"
    $ws.Range("E26").Value = '  -3.61%  '
    $ws.Range("D27").NumberFormat = "@"
    $ws.Range("D27").Value = '154.07'
    $ws.Range("D27").Style = "Normal"
    $ws.Range("E27").Value = '  -0.15%  '
    $ws.Range("D28").NumberFormat = "@"
    $ws.Range("D28").Value = '19.89'
    $ws.Range("D28").Style = "Normal"
    $ws.Range("E28").Value = '  -0.94%  '
    $ws.Range("D29").NumberFormat = "@"
    $ws.Range("D29").Value = '2.051'
    $ws.Range("D29").Style = "Normal"
    $ws.Range("E29").Value = '  -4.56%  '
    $ws.Range("D30").NumberFormat = "@"
    $ws.Range("D30").Value = '5.434'
    $ws.Range("D30").Style = "Normal"
    $ws.Range("E30").Value = '  -5.56%  '
    $ws.Range("D31").NumberFormat = "@"
    $ws.Range("D31").Value = '121.19'
    $ws.Range("D31").Style = "Normal"
    $ws.Range("E31").Value = '  -2.10%  '
    $ws.Range("D32").NumberFormat = "@"
    $ws.Range("D32").Value = '0.09311'
    $ws.Range("D32").Style = "Normal"
    $ws.Range("E32").Value = '  -3.19%  '
    $ws.Range("D33").NumberFormat = "@"
    $ws.Range("D33").Value = '0.9374'
    $ws.Range("D33").Style = "Normal"
    $ws.Range("E33").Value = '  -4.16%  '
    $ws.Range("D34").NumberFormat = "@"
    $ws.Range("D34").Value = '1.445'
    $ws.Range("D34").Style = "Normal"
    $ws.Range("E34").Value = '  -0.62%  '
    $ws.Range("D35").NumberFormat = "@"
    $ws.Range("D35").Value = '3.585'
    $ws.Range("D35").Style = "Normal"
    $ws.Range("E35").Value = '  -1.43%  '
    $ws.Range("D36").NumberFormat = "@"
    $ws.Range("D36").Value = '5.248'
    $ws.Range("D36").Style = "Normal"
    $ws.Range("E36").Value = '  -5.77%  '
    $ws.Range("E37").Value = '  -2.42%  '
    $ws.Range("D38").NumberFormat = "@"
    $ws.Range("D38").Value = '0.02222'
    $ws.Range("D38").Style = "Normal"
    $ws.Range("E38").Value = '  -4.10%  '
    $ws.Range("D39").NumberFormat = "@"
    $ws.Range("D39").Value = '1.204'
    $ws.Range("D39").Style = "Normal"
    $ws.Range("E39").Value = '  -2.60%  '
    $ws.Range("D40").NumberFormat = "@"
    $ws.Range("D40").Value = '8.083'
    $ws.Range("D40").Style = "Normal"
    $ws.Range("E40").Value = '  -10.05%  '
    $ws.Range("D42").NumberFormat = "@"
    $ws.Range("D42").Value = '0.5905'
    $ws.Range("D42").Style = "Normal"
    $ws.Range("E42").Value = '  -4.26%  '
    $ws.Range("D43").NumberFormat = "@"
    $ws.Range("D43").Value = '0.1882'
    $ws.Range("D43").Style = "Normal"
    $ws.Range("E43").Value = '  -1.23%  '
    $ws.Range("D44").NumberFormat = "@"
    $ws.Range("D44").Value = '10.14'
    $ws.Range("D44").Style = "Normal"
    $ws.Range("E44").Value = '  -8.44%  '
    $ws.Range("D45").NumberFormat = "@"
    $ws.Range("D45").Value = '1.267'
    $ws.Range("D45").Style = "Normal"
    $ws.Range("E45").Value = '  -3.07%  '
    $ws.Range("D46").NumberFormat = "@"
    $ws.Range("D46").Value = '0.5586'
    $ws.Range("D46").Style = "Normal"
    $ws.Range("E46").Value = '  -5.10%  '
    $ws.Range("D47").NumberFormat = "@"
    $ws.Range("D47").Value = '12.01'
    $ws.Range("D47").Style = "Normal"
    $ws.Range("E47").Value = '  -6.03%  '
    $ws.Range("D48").NumberFormat = "@"
    $ws.Range("D48").Value = '3.373'
    $ws.Range("D48").Style = "Normal"
    $ws.Range("E48").Value = '  -2.88%  '
    $ws.Range("D49").NumberFormat = "@"
    $ws.Range("D49").Value = '1.909'
    $ws.Range("D49").Style = "Normal"
    $ws.Range("E49").Value = '  -6.34%  '
    $ws.Range("D50").NumberFormat = "@"
    $ws.Range("D50").Value = '0.06729'
    $ws.Range("D50").Style = "Normal"
    $ws.Range("E50").Value = '  -1.61%  '
    $ws.Range("D51").NumberFormat = "@"
    $ws.Range("D51").Value = '107.82'
    $ws.Range("D51").Style = "Normal"
    $ws.Range("E51").Value = '  -2.24%  '
